# Region 4B.xlsx edit:
#  - Fix sheet name typo: "Puerto Prinsesa" -> "Puerto Princesa"
#  - Update the saved cell selection on that sheet from D14 to D18
#    (user had scrolled/clicked down a few rows before saving)

$wb = $excel.ActiveWorkbook

# Rename the second sheet (tab order is Calapan, Puerto Prinsesa)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "Puerto Princesa"

# Make it the active sheet and move the selection to D18, matching
# the saved <selection activeCell="D18" sqref="D18"/> in the sheet XML.
$ws.Activate()
$ws.Range("D18").Select()
